$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (font/style) from the last existing row as a template for the new rows
$ws.Range("A1021:G1021").Copy()
$ws.Range("A1022:G1051").PasteSpecial(-4122)

# --- Pass 1: fill in the ID column (A) for all 30 new rows first (Unit 35, U35_01..U35_30) ---
$ws.Range("A1022").Value = 'U35_01'
$ws.Range("A1023").Value = 'U35_02'
$ws.Range("A1024").Value = 'U35_03'
$ws.Range("A1025").Value = 'U35_04'
$ws.Range("A1026").Value = 'U35_05'
$ws.Range("A1027").Value = 'U35_06'
$ws.Range("A1028").Value = 'U35_07'
$ws.Range("A1029").Value = 'U35_08'
$ws.Range("A1030").Value = 'U35_09'
$ws.Range("A1031").Value = 'U35_10'
$ws.Range("A1032").Value = 'U35_11'
$ws.Range("A1033").Value = 'U35_12'
$ws.Range("A1034").Value = 'U35_13'
$ws.Range("A1035").Value = 'U35_14'
$ws.Range("A1036").Value = 'U35_15'
$ws.Range("A1037").Value = 'U35_16'
$ws.Range("A1038").Value = 'U35_17'
$ws.Range("A1039").Value = 'U35_18'
$ws.Range("A1040").Value = 'U35_19'
$ws.Range("A1041").Value = 'U35_20'
$ws.Range("A1042").Value = 'U35_21'
$ws.Range("A1043").Value = 'U35_22'
$ws.Range("A1044").Value = 'U35_23'
$ws.Range("A1045").Value = 'U35_24'
$ws.Range("A1046").Value = 'U35_25'
$ws.Range("A1047").Value = 'U35_26'
$ws.Range("A1048").Value = 'U35_27'
$ws.Range("A1049").Value = 'U35_28'
$ws.Range("A1050").Value = 'U35_29'
$ws.Range("A1051").Value = 'U35_30'

# --- Pass 2: fill in the rest of each row (UnitID, Word, Answer, Sentence, Phrase, Part-of-speech) ---
$ws.Range("B1022").Value = 35
$ws.Range("C1022").Value = 'Kỷ niệm'
$ws.Range("D1022").Value = 'Anniversary'
$ws.Range("E1022").Value = 'We are celebrating our wedding anniversary'
$ws.Range("F1022").Value = 'A wedding anniversary'
$ws.Range("G1022").Value = 'N'
$ws.Range("B1023").Value = 35
$ws.Range("C1023").Value = 'Cuộc hôn nhân'
$ws.Range("D1023").Value = 'Marriage'
$ws.Range("E1023").Value = 'There are arranged marriages in india'
$ws.Range("F1023").Value = 'An arranged marriage / hôn nhân sắp đặt'
$ws.Range("G1023").Value = 'N'
$ws.Range("B1024").Value = 35
$ws.Range("C1024").Value = 'Lãng mạn'
$ws.Range("D1024").Value = 'Romantic'
$ws.Range("E1024").Value = 'The couple has a beautiful romantic story'
$ws.Range("F1024").Value = 'A romantic story'
$ws.Range("G1024").Value = 'Adj'
$ws.Range("B1025").Value = 35
$ws.Range("C1025").Value = 'Đỏ mặt'
$ws.Range("D1025").Value = 'Blush'
$ws.Range("E1025").Value = 'His action made me blush with embarrassment'
$ws.Range("F1025").Value = 'blush with embarrassment / đỏ mặt ngượng ngùng'
$ws.Range("G1025").Value = 'V'
$ws.Range("B1026").Value = 35
$ws.Range("C1026").Value = 'Sự kết nối'
$ws.Range("D1026").Value = 'Connection'
$ws.Range("E1026").Value = 'The connection between them is strong'
$ws.Range("F1026").Value = 'a connection between something or somebody'
$ws.Range("G1026").Value = 'N'
$ws.Range("B1027").Value = 35
$ws.Range("C1027").Value = 'Ngay lập tức'
$ws.Range("D1027").Value = 'Instantly'
$ws.Range("E1027").Value = 'Her voice is instantly recognizable.'
$ws.Range("F1027").Value = 'instantly recognizable.'
$ws.Range("G1027").Value = 'Adv'
$ws.Range("B1028").Value = 35
$ws.Range("C1028").Value = 'Rủ đi chơi'
$ws.Range("D1028").Value = 'Ask out'
$ws.Range("E1028").Value = 'Can I ask you out on a date?'
$ws.Range("F1028").Value = 'Ask somebody out on a date / mời ai một buổi hẹn'
$ws.Range("G1028").Value = 'V'
$ws.Range("B1029").Value = 35
$ws.Range("C1029").Value = 'Vững chắc'
$ws.Range("D1029").Value = 'Steady'
$ws.Range("E1029").Value = 'They try to maintain a steady relationship'
$ws.Range("F1029").Value = 'steady relationship / mối quan hệ vững chắc'
$ws.Range("G1029").Value = 'Adj'
$ws.Range("B1030").Value = 35
$ws.Range("C1030").Value = 'Phải lòng'
$ws.Range("D1030").Value = 'Fall for'
$ws.Range("E1030").Value = 'Don''t make me fall for you.'
$ws.Range("F1030").Value = 'fall for somebody'
$ws.Range("G1030").Value = 'V'
$ws.Range("B1031").Value = 35
$ws.Range("C1031").Value = 'Hấp dẫn, quyến rũ'
$ws.Range("D1031").Value = 'Attractive'
$ws.Range("E1031").Value = 'She is an attractive person.'
$ws.Range("F1031").Value = 'an attractive person / một người hấp dẫn'
$ws.Range("G1031").Value = 'Adj'
$ws.Range("B1032").Value = 35
$ws.Range("C1032").Value = 'Rõ ràng'
$ws.Range("D1032").Value = 'Obviously'
$ws.Range("E1032").Value = 'Diet and exercise are obviously important.'
$ws.Range("F1032").Value = 'to be obviously important'
$ws.Range("G1032").Value = 'Adv'
$ws.Range("B1033").Value = 35
$ws.Range("C1033").Value = 'Chia tay'
$ws.Range("D1033").Value = 'Break up'
$ws.Range("E1033").Value = 'She cries when the actor breaks up with his girlfriend'
$ws.Range("F1033").Value = 'break up with somebody'
$ws.Range("G1033").Value = 'V'
$ws.Range("B1034").Value = 35
$ws.Range("C1034").Value = 'Liên hệ'
$ws.Range("D1034").Value = 'Contact'
$ws.Range("E1034").Value = 'We will contact you about the decision.'
$ws.Range("F1034").Value = 'contact someone'
$ws.Range("G1034").Value = 'V'
$ws.Range("B1035").Value = 35
$ws.Range("C1035").Value = 'Dịu dàng, nhẹ nhàng'
$ws.Range("D1035").Value = 'Gentle'
$ws.Range("E1035").Value = 'The singer has such a gentle voice.'
$ws.Range("F1035").Value = 'a gentle voice / một chất giọng nhẹ nhàng'
$ws.Range("G1035").Value = 'Adj'
$ws.Range("B1036").Value = 35
$ws.Range("C1036").Value = 'Lời nối dối'
$ws.Range("D1036").Value = 'Lie'
$ws.Range("E1036").Value = 'I don''t think it''s okay to tell a white lie.'
$ws.Range("F1036").Value = 'a white lie / lời nói dối vô hại'
$ws.Range("G1036").Value = 'N'
$ws.Range("B1037").Value = 35
$ws.Range("C1037").Value = 'Tình huống'
$ws.Range("D1037").Value = 'Situation'
$ws.Range("E1037").Value = 'I am facing a difficult situation now'
$ws.Range("F1037").Value = 'a difficult situation'
$ws.Range("G1037").Value = 'N'
$ws.Range("B1038").Value = 35
$ws.Range("C1038").Value = 'Lời bào chữa'
$ws.Range("D1038").Value = 'Excuse'
$ws.Range("E1038").Value = 'There is no excuse for arriving late.'
$ws.Range("F1038").Value = 'an excuse for doing something / một cái cớ để làm một cái gì đó'
$ws.Range("G1038").Value = 'N'
$ws.Range("B1039").Value = 35
$ws.Range("C1039").Value = 'Thư'
$ws.Range("D1039").Value = 'Letter'
$ws.Range("E1039").Value = 'I am sending a letter of complaint soon.'
$ws.Range("F1039").Value = 'a letter of complaint / thư khiếu nại'
$ws.Range("G1039").Value = 'N'
$ws.Range("B1040").Value = 35
$ws.Range("C1040").Value = 'Ý định'
$ws.Range("D1040").Value = 'Intention'
$ws.Range("E1040").Value = 'My intention of borrowing your car is to impress her'
$ws.Range("F1040").Value = 'intention of doing something'
$ws.Range("G1040").Value = 'N'
$ws.Range("B1041").Value = 35
$ws.Range("C1041").Value = 'Tha thứ'
$ws.Range("D1041").Value = 'Forgive'
$ws.Range("E1041").Value = 'Would you ever forgive me?'
$ws.Range("F1041").Value = 'forgive somebody / tha thứ cho ai đó'
$ws.Range("G1041").Value = 'V'
$ws.Range("B1042").Value = 35
$ws.Range("C1042").Value = 'Cô dâu'
$ws.Range("D1042").Value = 'Bride'
$ws.Range("E1042").Value = 'He introduces his new bride.'
$ws.Range("F1042").Value = 'a new bride / vợ mới cưới'
$ws.Range("G1042").Value = 'N'
$ws.Range("B1043").Value = 35
$ws.Range("C1043").Value = 'Chú rể'
$ws.Range("D1043").Value = 'Groom'
$ws.Range("E1043").Value = 'Let us toast to the (nâng ly chúc mừng) bride and groom.'
$ws.Range("F1043").Value = 'The bride and groom / cô dâu chú rể'
$ws.Range("G1043").Value = 'N'
$ws.Range("B1044").Value = 35
$ws.Range("C1044").Value = 'Hôn lễ'
$ws.Range("D1044").Value = 'Wedding'
$ws.Range("E1044").Value = 'Here is a wedding present for the bride.'
$ws.Range("F1044").Value = 'a wedding present / một món quà cưới'
$ws.Range("G1044").Value = 'N'
$ws.Range("B1045").Value = 35
$ws.Range("C1045").Value = 'Cặp đôi'
$ws.Range("D1045").Value = 'Couple'
$ws.Range("E1045").Value = 'We are officially a married couple.'
$ws.Range("F1045").Value = 'a married couple / một cặp vợ chồng'
$ws.Range("G1045").Value = 'N'
$ws.Range("B1046").Value = 35
$ws.Range("C1046").Value = 'Cảm xúc'
$ws.Range("D1046").Value = 'Feeling'
$ws.Range("F1046").Value = 'a feeling of sadness / cảm giác buồn bã'
$ws.Range("E1046").Value = 'Losing (thất bại) gives me a feeling of sadness'
$ws.Range("G1046").Value = 'N'
$ws.Range("B1047").Value = 35
$ws.Range("C1047").Value = 'Ấm áp, nồng nhiệt'
$ws.Range("D1047").Value = 'Warm'
$ws.Range("E1047").Value = 'The host extends a warm welcome to us.'
$ws.Range("F1047").Value = 'a warm welcome / một sự chào đón nồng nhiệt'
$ws.Range("G1047").Value = 'Adj'
$ws.Range("B1048").Value = 35
$ws.Range("C1048").Value = 'Đối diện'
$ws.Range("D1048").Value = 'Opposite'
$ws.Range("E1048").Value = 'he sits opposite you in class.'
$ws.Range("F1048").Value = 'to be opposite something or someone'
$ws.Range("G1048").Value = 'Adj'
$ws.Range("B1049").Value = 35
$ws.Range("C1049").Value = 'Kỳ cục'
$ws.Range("D1049").Value = 'Weird'
$ws.Range("E1049").Value = 'After lunch, my stomach had a weird feeling.'
$ws.Range("F1049").Value = 'a weird feeling / một cảm giác kỳ lạ'
$ws.Range("G1049").Value = 'Adj'
$ws.Range("B1050").Value = 35
$ws.Range("C1050").Value = 'Sẵn sàng'
$ws.Range("D1050").Value = 'Willing'
$ws.Range("E1050").Value = 'Would you be willing to help me?'
$ws.Range("F1050").Value = 'To be willing to do something'
$ws.Range("G1050").Value = 'Adj'
$ws.Range("B1051").Value = 35
$ws.Range("C1051").Value = 'Khoảnh khắc'
$ws.Range("D1051").Value = 'Moment'
$ws.Range("E1051").Value = 'She glanced at me for a brief moment'
$ws.Range("F1051").Value = 'a brief moment / một khoảnh khắc ngắn ngủi'
$ws.Range("G1051").Value = 'N'

# Row heights to match the rest of the sheet (17pt, driven by the 13pt Times New Roman font)
$ws.Range("A1022:G1051").RowHeight = 17

# Column A needs to widen slightly to fit the new content
$ws.Columns.Item(1).ColumnWidth = 7.666666666666667

# Leave the cursor where the author would continue entering the next word (row 1052)
[void]$ws.Range("C1052").Select()
